# Commit: "Added Run Mode concept"
# Adds a new "Test_Cases" worksheet (with a Sr.No / TestCaseName / Run_Mode
# table) as the first sheet in the workbook, ahead of the existing
# "TestData" and "NewSheet" sheets.

$wb = $excel.ActiveWorkbook

# "TestData" now shows a different selected cell (B17 instead of H17) in
# the target workbook, so switch to it and move the selection there first.
$testData = $wb.Worksheets.Item("TestData")
$testData.Activate()
$testData.Range("B17").Select()

# Worksheets.Add() with no args inserts a new sheet immediately before the
# currently active sheet - i.e. right at the front, matching the target
# workbook order: Test_Cases, TestData, NewSheet.
$ws = $wb.Worksheets.Add()
$ws.Name = "Test_Cases"

# ---- Header row -----------------------------------------------------
$ws.Range("A1").Value = "Sr.No"
$ws.Range("B1").Value = "TestCaseName"
$ws.Range("C1").Value = "Run_Mode"

$header = $ws.Range("A1:C1")
$header.Font.Bold = $true
$header.Interior.Color = 65535
$header.Borders.LineStyle = 1
$header.HorizontalAlignment = -4108

# ---- Data rows --------------------------------------------------------
$data = @(
    @(1, "CreateAnNewAccount1", "Y"),
    @(2, "CreateAnNewAccount2", "Y"),
    @(3, "CreateAnNewAccount3", "N"),
    @(4, "TC04_Create_an_Account4", "N"),
    @(5, "TC05_Create_an_Account5", "N")
)

$r = 2
foreach ($row in $data) {
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $r = $r + 1
}

$body = $ws.Range("A2:C6")
$body.Borders.LineStyle = 1
$body.HorizontalAlignment = -4108

# ---- Column widths ------------------------------------------------
$ws.Range("A1").ColumnWidth = 5.36328125
$ws.Range("B1").ColumnWidth = 33.90625
$ws.Range("C1").ColumnWidth = 10

# ---- Sheet view housekeeping ---------------------------------------
# Make Test_Cases the active tab, with D17 remembered as its selection.
$ws.Activate()
$ws.Range("D17").Select()
